$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178-247 down to 179-248.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with a new data record (same shape as the
# surrounding "Feria Lagunitas de Puerto Montt" / Apio / Primera rows), carrying
# the style of column D (date) that the insert already copied down.
$ws.Cells.Item(178, 1).Value = 4
$ws.Cells.Item(178, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(178, 3).Value = "Los Lagos"
$ws.Cells.Item(178, 4).Value = 44704
$ws.Cells.Item(178, 5).Value = 10
$ws.Cells.Item(178, 6).Value = 100112017
$ws.Cells.Item(178, 7).Value = "Apio"
$ws.Cells.Item(178, 8).Value = "Americana (o)"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 20
$ws.Cells.Item(178, 11).Value = 12000
$ws.Cells.Item(178, 12).Value = 12000
$ws.Cells.Item(178, 13).Value = 12000
$ws.Cells.Item(178, 14).Value = "`$/docena de matas"
$ws.Cells.Item(178, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(178, 16).Value = 2000
$ws.Cells.Item(178, 17).Value = 6
$ws.Cells.Item(178, 18).Value = "Hortaliza"
